# Apply data updates to the "Quarterly Data Levels" sheet
# (columns M = Savings Rate, O = Employment level) per the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Quarterly Data Levels")

$updates = @{
    "M71" = 14
    "M72" = 14.7
    "M73" = 14.6
    "M74" = 13.6
    "M75" = 13.7
    "M76" = 12.8
    "M77" = 13.3
    "M78" = 13.1
    "M79" = 11.1
    "M80" = 11
    "M81" = 10.9
    "M82" = 13.1
    "M83" = 14
    "M84" = 12.6
    "M85" = 11.5
    "M86" = 10.7
    "M87" = 12.3
    "M88" = 14.4
    "M89" = 12.6
    "M90" = 14.3
    "M91" = 15
    "M92" = 15.1
    "M93" = 14.1
    "M97" = 12.7
    "M98" = 14.1
    "M99" = 13.5
    "M100" = 13.4
    "M101" = 14.4
    "M102" = 25.1
    "M103" = 33.4
    "M104" = 20.7
    "M105" = 24.5
    "M106" = 26.2
    "M107" = 23.1
    "M108" = 19.6
    "M109" = 17.9
    "M110" = 15.9
    "M111" = 15.7
    "M112" = 14.8
    "M113" = 12.7
    "O113" = 2399700
    "M114" = 13
    "M115" = 11.4
    "M117" = 10.8
    "M118" = 14.7
    "O121" = 2519733.33333333
    "M122" = 13.2
    "M123" = 12.5
    "O123" = 2553533.33333333
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
